# Append a new customer row (phone 51616172) with total_points reset to 0.
# Mirrors: Update points 51616172 -> 0.00

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Cells.Item(63, 1).Row

# Column A (phone): stored as text, not a number, so the leading zeros /
# exact digit string are preserved (leading apostrophe forces text entry,
# same as typing '51616172 into the cell).
$ws.Cells.Item($newRow, 1).Value = "'51616172"

# Column B (birthday): unknown/blank for this customer - store as an
# empty text value (matches the existing blank-birthday rows in the sheet).
$ws.Cells.Item($newRow, 2).Value = "'"

# Column C (total_points): reset to 0.
$ws.Cells.Item($newRow, 3).Value = 0
